$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 243
$ws.Range("I2").Value = 243
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 243
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -130
$ws.Range("N2").ClearContents()

# ALC row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 821.6
$ws.Range("I96").Value = 976.5
$ws.Range("J96").Value = 589.25
$ws.Range("K96").Value = 2929.5
$ws.Range("L96").Value = 1767.75
$ws.Range("M96").Value = -1556.5
$ws.Range("N96").Value = -4513.75

# ALC row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 47355.957
$ws.Range("I132").Value = 56510.79
$ws.Range("J132").Value = 3870.5
$ws.Range("K132").Value = 169532.37
$ws.Range("L132").Value = 11611.5
$ws.Range("M132").Value = -167002.37
$ws.Range("N132").Value = -16671.5

# ALC row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 16668607
$ws.Range("I137").Value = 40909750
$ws.Range("K137").Value = 122729250
$ws.Range("M137").Value = -122726700

# ALC row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 4649.5557
$ws.Range("I138").Value = 4974.8335
$ws.Range("J138").Value = 3999
$ws.Range("K138").Value = 14924.5005
$ws.Range("L138").Value = 11997
$ws.Range("M138").Value = -9784.500499999998
$ws.Range("N138").Value = -22277

$ws = $wb.Worksheets.Item("ARM")
# ARM row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()

# ARM row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 2499.95
$ws.Range("I61").Value = 2354.7778
$ws.Range("J61").Value = 3806.5
$ws.Range("K61").Value = 2354.7778
$ws.Range("L61").Value = 3806.5
$ws.Range("M61").Value = -2142.7778
$ws.Range("N61").Value = -4230.5

# ARM row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 559557.75
$ws.Range("I74").Value = 3953.8372
$ws.Range("J74").Value = 3972553
$ws.Range("K74").Value = 3953.8372
$ws.Range("L74").Value = 3972553
$ws.Range("M74").Value = -3079.8372
$ws.Range("N74").Value = -3974301

# ARM row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 559557.75
$ws.Range("I77").Value = 3953.8372
$ws.Range("J77").Value = 3972553
$ws.Range("K77").Value = 19769.186
$ws.Range("L77").Value = 19862765
$ws.Range("M77").Value = -15401.186
$ws.Range("N77").Value = -19871501

# ARM row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 45458304
$ws.Range("I132").Value = 2854.625
$ws.Range("J132").Value = 166672830
$ws.Range("K132").Value = 8563.875
$ws.Range("L132").Value = 500018490
$ws.Range("M132").Value = -6033.875
$ws.Range("N132").Value = -500023550

# ARM row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 2499.95
$ws.Range("I136").Value = 2354.7778
$ws.Range("J136").Value = 3806.5
$ws.Range("K136").Value = 7064.3334
$ws.Range("L136").Value = 11419.5
$ws.Range("M136").Value = -4514.3334
$ws.Range("N136").Value = -16519.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 81 (Leve Item ID 42300)
$ws.Range("H81").Value = 19592.666
$ws.Range("J81").Value = 19592.666
$ws.Range("L81").Value = 19592.666
$ws.Range("N81").Value = -21714.666

# BSM row 84 (Leve Item ID 42300)
$ws.Range("H84").Value = 19592.666
$ws.Range("J84").Value = 19592.666
$ws.Range("L84").Value = 58777.99800000001
$ws.Range("N84").Value = -69385.99800000001

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4625.66
$ws.Range("I31").Value = 2862.1333
$ws.Range("J31").Value = 5321.7896
$ws.Range("K31").Value = 2862.1333
$ws.Range("L31").Value = 5321.7896
$ws.Range("M31").Value = -2567.1333
$ws.Range("N31").Value = -5911.7896

# CRP row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4625.66
$ws.Range("I34").Value = 2862.1333
$ws.Range("J34").Value = 5321.7896
$ws.Range("K34").Value = 2862.1333
$ws.Range("L34").Value = 5321.7896
$ws.Range("M34").Value = -2660.1333
$ws.Range("N34").Value = -5725.7896

# CRP row 70 (Leve Item ID 12011)
$ws.Range("H70").Value = 28835.8
$ws.Range("J70").Value = 28835.8
$ws.Range("L70").Value = 28835.8
$ws.Range("N70").Value = -29465.8

# CRP row 73 (Leve Item ID 12011)
$ws.Range("H73").Value = 28835.8
$ws.Range("J73").Value = 28835.8
$ws.Range("L73").Value = 28835.8
$ws.Range("N73").Value = -31019.8

# CRP row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 3011.4285
$ws.Range("I134").Value = 2767.25
$ws.Range("J134").Value = 3337
$ws.Range("K134").Value = 8301.75
$ws.Range("L134").Value = 10011
$ws.Range("M134").Value = -5766.75
$ws.Range("N134").Value = -15081

$ws = $wb.Worksheets.Item("CUL")
# CUL row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 3956.05
$ws.Range("I68").Value = 1810.8
$ws.Range("J68").Value = 4671.1333
$ws.Range("K68").Value = 5432.4
$ws.Range("L68").Value = 14013.3999
$ws.Range("M68").Value = -4621.4
$ws.Range("N68").Value = -15635.3999

# CUL row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 3956.05
$ws.Range("I71").Value = 1810.8
$ws.Range("J71").Value = 4671.1333
$ws.Range("K71").Value = 16297.2
$ws.Range("L71").Value = 42040.1997
$ws.Range("M71").Value = -12241.2
$ws.Range("N71").Value = -50152.1997

# CUL row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 323.85715
$ws.Range("J107").Value = 289.5
$ws.Range("L107").Value = 868.5
$ws.Range("N107").Value = -4708.5

# CUL row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1238
$ws.Range("I132").Value = 1214.8
$ws.Range("J132").Value = 1315.3334
$ws.Range("K132").Value = 10933.2
$ws.Range("L132").Value = 11838.0006
$ws.Range("M132").Value = -8403.199999999999
$ws.Range("N132").Value = -16898.0006

$ws = $wb.Worksheets.Item("GSM")
# GSM row 75 (Leve Item ID 11008)
$ws.Range("H75").Value = 43348
$ws.Range("J75").Value = 43348
$ws.Range("L75").Value = 43348
$ws.Range("N75").Value = -45096

# GSM row 78 (Leve Item ID 11008)
$ws.Range("H78").Value = 43348
$ws.Range("J78").Value = 43348
$ws.Range("L78").Value = 130044
$ws.Range("N78").Value = -138780

# GSM row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 7408.636
$ws.Range("I132").Value = 7816.5
$ws.Range("K132").Value = 23449.5
$ws.Range("M132").Value = -20919.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 74 (Leve Item ID 19022)
$ws.Range("H74").Value = 118621.29
$ws.Range("I74").Value = 14000
$ws.Range("J74").Value = 136058.17
$ws.Range("K74").Value = 14000
$ws.Range("L74").Value = 136058.17
$ws.Range("M74").Value = -13064
$ws.Range("N74").Value = -137930.17

# WVR row 77 (Leve Item ID 19022)
$ws.Range("H77").Value = 118621.29
$ws.Range("I77").Value = 14000
$ws.Range("J77").Value = 136058.17
$ws.Range("K77").Value = 42000
$ws.Range("L77").Value = 408174.51
$ws.Range("M77").Value = -37320
$ws.Range("N77").Value = -417534.51

